# Update "想去人数" (attendance interest count) values on the 展览 and 全部类型
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3 = 2337
    4 = 416
    5 = 82
    6 = 6475
    7 = 326
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
